$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.820425
$ws.Range("H2").Value = 11.461275
$ws.Range("I2").Value = 0.02049663039797357
$ws.Range("J2").Value = 0.02049663039797357
$ws.Range("M2").Value = 86.89540866666668
$ws.Range("N2").Value = 260.686226
$ws.Range("O2").Value = 0.319779657009892
$ws.Range("P2").Value = 0.3197796570098919
$ws.Range("Q2").Value = 331.9773916553501
$ws.Range("R2").Value = 2987.79652489815
$ws.Range("S2").Value = 0.006554405438522513
$ws.Range("T2").Value = 0.006554405438522511
$ws.Range("G3").Value = 3.820425
$ws.Range("H3").Value = 11.461275
$ws.Range("I3").Value = 0.02049663039797357
$ws.Range("J3").Value = 0.02049663039797357
$ws.Range("O3").Value = 0.1999969065479545
$ws.Range("P3").Value = 0.1999969065479545
$ws.Range("Q3").Value = 207.625625706625
$ws.Range("R3").Value = 1868.630631359625
$ws.Range("S3").Value = 0.004099262674251484
$ws.Range("T3").Value = 0.004099262674251483
$ws.Range("G4").Value = 3.820425
$ws.Range("H4").Value = 11.461275
$ws.Range("I4").Value = 0.02049663039797357
$ws.Range("J4").Value = 0.02049663039797357
$ws.Range("M4").Value = 60.92601633333334
$ws.Range("N4").Value = 182.778049
$ws.Range("O4").Value = 0.224210932487692
$ws.Range("P4").Value = 0.224210932487692
$ws.Range("Q4").Value = 232.763275950275
$ws.Range("R4").Value = 2094.869483552475
$ws.Range("S4").Value = 0.004595568614385227
$ws.Range("T4").Value = 0.004595568614385226
$ws.Range("G5").Value = 3.820425
$ws.Range("H5").Value = 11.461275
$ws.Range("I5").Value = 0.02049663039797357
$ws.Range("J5").Value = 0.02049663039797357
$ws.Range("M5").Value = 7.809668333333332
$ws.Range("N5").Value = 23.429005
$ws.Range("O5").Value = 0.02873998867505581
$ws.Range("P5").Value = 0.02873998867505581
$ws.Range("Q5").Value = 29.836252142375
$ws.Range("R5").Value = 268.526269281375
$ws.Range("S5").Value = 0.000589072925514565
$ws.Range("T5").Value = 0.000589072925514565
$ws.Range("G6").Value = 3.820425
$ws.Range("H6").Value = 11.461275
$ws.Range("I6").Value = 0.02049663039797357
$ws.Range("J6").Value = 0.02049663039797357
$ws.Range("M6").Value = 61.75795633333333
$ws.Range("N6").Value = 185.273869
$ws.Range("O6").Value = 0.2272725152794058
$ws.Range("P6").Value = 0.2272725152794058
$ws.Range("Q6").Value = 235.941640324775
$ws.Range("R6").Value = 2123.474762922975
$ws.Range("S6").Value = 0.004658320745299782
$ws.Range("T6").Value = 0.004658320745299781
$ws.Range("I7").Value = 0.7542622677884155
$ws.Range("J7").Value = 0.7542622677884157
$ws.Range("M7").Value = 86.89540866666668
$ws.Range("N7").Value = 260.686226
$ws.Range("O7").Value = 0.319779657009892
$ws.Range("P7").Value = 0.3197796570098919
$ws.Range("Q7").Value = 12216.54561860097
$ws.Range("R7").Value = 109948.9105674087
$ws.Range("S7").Value = 0.2411977292888828
$ws.Range("T7").Value = 0.2411977292888828
$ws.Range("I8").Value = 0.7542622677884155
$ws.Range("J8").Value = 0.7542622677884157
$ws.Range("O8").Value = 0.1999969065479545
$ws.Range("P8").Value = 0.1999969065479545
$ws.Range("S8").Value = 0.150850120283528
$ws.Range("T8").Value = 0.150850120283528
$ws.Range("I9").Value = 0.7542622677884155
$ws.Range("J9").Value = 0.7542622677884157
$ws.Range("M9").Value = 60.92601633333334
$ws.Range("N9").Value = 182.778049
$ws.Range("O9").Value = 0.224210932487692
$ws.Range("P9").Value = 0.224210932487692
$ws.Range("Q9").Value = 8565.532624985652
$ws.Range("R9").Value = 77089.79362487086
$ws.Range("S9").Value = 0.1691138464011219
$ws.Range("T9").Value = 0.1691138464011219
$ws.Range("I10").Value = 0.7542622677884155
$ws.Range("J10").Value = 0.7542622677884157
$ws.Range("M10").Value = 7.809668333333332
$ws.Range("N10").Value = 23.429005
$ws.Range("O10").Value = 0.02873998867505581
$ws.Range("P10").Value = 0.02873998867505581
$ws.Range("Q10").Value = 1097.954091294912
$ws.Range("R10").Value = 9881.586821654209
$ws.Range("S10").Value = 0.02167748903426098
$ws.Range("T10").Value = 0.02167748903426098
$ws.Range("I11").Value = 0.7542622677884155
$ws.Range("J11").Value = 0.7542622677884157
$ws.Range("M11").Value = 61.75795633333333
$ws.Range("N11").Value = 185.273869
$ws.Range("O11").Value = 0.2272725152794058
$ws.Range("P11").Value = 0.2272725152794058
$ws.Range("Q11").Value = 8682.494304755477
$ws.Range("R11").Value = 78142.4487427993
$ws.Range("S11").Value = 0.171423082780622
$ws.Range("T11").Value = 0.171423082780622
$ws.Range("G12").Value = 30.51067
$ws.Range("H12").Value = 91.53201
$ws.Range("I12").Value = 0.1636901460399144
$ws.Range("J12").Value = 0.1636901460399144
$ws.Range("M12").Value = 86.89540866666668
$ws.Range("N12").Value = 260.686226
$ws.Range("O12").Value = 0.319779657009892
$ws.Range("P12").Value = 0.3197796570098919
$ws.Range("Q12").Value = 2651.237138343807
$ws.Range("R12").Value = 23861.13424509426
$ws.Range("S12").Value = 0.05234477875654297
$ws.Range("T12").Value = 0.05234477875654295
$ws.Range("G13").Value = 30.51067
$ws.Range("H13").Value = 91.53201
$ws.Range("I13").Value = 0.1636901460399144
$ws.Range("J13").Value = 0.1636901460399144
$ws.Range("O13").Value = 0.1999969065479545
$ws.Range("P13").Value = 0.1999969065479545
$ws.Range("Q13").Value = 1658.139329911817
$ws.Range("R13").Value = 14923.25396920635
$ws.Range("S13").Value = 0.0327375228403658
$ws.Range("T13").Value = 0.0327375228403658
$ws.Range("G14").Value = 30.51067
$ws.Range("H14").Value = 91.53201
$ws.Range("I14").Value = 0.1636901460399144
$ws.Range("J14").Value = 0.1636901460399144
$ws.Range("M14").Value = 60.92601633333334
$ws.Range("N14").Value = 182.778049
$ws.Range("O14").Value = 0.224210932487692
$ws.Range("P14").Value = 0.224210932487692
$ws.Range("Q14").Value = 1858.893578760943
$ws.Range("R14").Value = 16730.04220884849
$ws.Range("S14").Value = 0.0367011202826557
$ws.Range("T14").Value = 0.03670112028265569
$ws.Range("G15").Value = 30.51067
$ws.Range("H15").Value = 91.53201
$ws.Range("I15").Value = 0.1636901460399144
$ws.Range("J15").Value = 0.1636901460399144
$ws.Range("M15").Value = 7.809668333333332
$ws.Range("N15").Value = 23.429005
$ws.Range("O15").Value = 0.02873998867505581
$ws.Range("P15").Value = 0.02873998867505581
$ws.Range("Q15").Value = 238.2782133277833
$ws.Range("R15").Value = 2144.50391995005
$ws.Range("S15").Value = 0.004704452943405373
$ws.Range("T15").Value = 0.004704452943405373
$ws.Range("G16").Value = 30.51067
$ws.Range("H16").Value = 91.53201
$ws.Range("I16").Value = 0.1636901460399144
$ws.Range("J16").Value = 0.1636901460399144
$ws.Range("M16").Value = 61.75795633333333
$ws.Range("N16").Value = 185.273869
$ws.Range("O16").Value = 0.2272725152794058
$ws.Range("P16").Value = 0.2272725152794058
$ws.Range("Q16").Value = 1884.276625560743
$ws.Range("R16").Value = 16958.48963004669
$ws.Range("S16").Value = 0.03720227121694463
$ws.Range("T16").Value = 0.03720227121694462
$ws.Range("G17").Value = 0.258813
$ws.Range("H17").Value = 0.776439
$ws.Range("I17").Value = 0.001388535150720334
$ws.Range("J17").Value = 0.001388535150720334
$ws.Range("M17").Value = 86.89540866666668
$ws.Range("N17").Value = 260.686226
$ws.Range("O17").Value = 0.319779657009892
$ws.Range("P17").Value = 0.3197796570098919
$ws.Range("Q17").Value = 22.48966140324601
$ws.Range("R17").Value = 202.406952629214
$ws.Range("S17").Value = 0.0004440252942435271
$ws.Range("T17").Value = 0.000444025294243527
$ws.Range("G18").Value = 0.258813
$ws.Range("H18").Value = 0.776439
$ws.Range("I18").Value = 0.001388535150720334
$ws.Range("J18").Value = 0.001388535150720334
$ws.Range("O18").Value = 0.1999969065479545
$ws.Range("P18").Value = 0.1999969065479545
$ws.Range("Q18").Value = 14.065506080085
$ws.Range("R18").Value = 126.589554720765
$ws.Range("S18").Value = 0.0002777027347771647
$ws.Range("T18").Value = 0.0002777027347771646
$ws.Range("G19").Value = 0.258813
$ws.Range("H19").Value = 0.776439
$ws.Range("I19").Value = 0.001388535150720334
$ws.Range("J19").Value = 0.001388535150720334
$ws.Range("M19").Value = 60.92601633333334
$ws.Range("N19").Value = 182.778049
$ws.Range("O19").Value = 0.224210932487692
$ws.Range("P19").Value = 0.224210932487692
$ws.Range("Q19").Value = 15.768445065279
$ws.Range("R19").Value = 141.916005587511
$ws.Range("S19").Value = 0.0003113247609349441
$ws.Range("T19").Value = 0.000311324760934944
$ws.Range("G20").Value = 0.258813
$ws.Range("H20").Value = 0.776439
$ws.Range("I20").Value = 0.001388535150720334
$ws.Range("J20").Value = 0.001388535150720334
$ws.Range("M20").Value = 7.809668333333332
$ws.Range("N20").Value = 23.429005
$ws.Range("O20").Value = 0.02873998867505581
$ws.Range("P20").Value = 0.02873998867505581
$ws.Range("Q20").Value = 2.021243690355
$ws.Range("R20").Value = 18.191193213195
$ws.Range("S20").Value = 0.00003990648450661932
$ws.Range("T20").Value = 0.00003990648450661931
$ws.Range("G21").Value = 0.258813
$ws.Range("H21").Value = 0.776439
$ws.Range("I21").Value = 0.001388535150720334
$ws.Range("J21").Value = 0.001388535150720334
$ws.Range("M21").Value = 61.75795633333333
$ws.Range("N21").Value = 185.273869
$ws.Range("O21").Value = 0.2272725152794058
$ws.Range("P21").Value = 0.2272725152794058
$ws.Range("Q21").Value = 15.983761952499
$ws.Range("R21").Value = 143.853857572491
$ws.Range("S21").Value = 0.0003155758762580792
$ws.Range("T21").Value = 0.0003155758762580792
$ws.Range("G22").Value = 11.213844
$ws.Range("H22").Value = 33.641532
$ws.Range("I22").Value = 0.0601624206229761
$ws.Range("J22").Value = 0.0601624206229761
$ws.Range("M22").Value = 86.89540866666668
$ws.Range("N22").Value = 260.686226
$ws.Range("O22").Value = 0.319779657009892
$ws.Range("P22").Value = 0.3197796570098919
$ws.Range("Q22").Value = 974.4315571042482
$ws.Range("R22").Value = 8769.884013938232
$ws.Range("S22").Value = 0.01923871823170015
$ws.Range("T22").Value = 0.01923871823170014
$ws.Range("G23").Value = 11.213844
$ws.Range("H23").Value = 33.641532
$ws.Range("I23").Value = 0.0601624206229761
$ws.Range("J23").Value = 0.0601624206229761
$ws.Range("O23").Value = 0.1999969065479545
$ws.Range("P23").Value = 0.1999969065479545
$ws.Range("Q23").Value = 609.42993962098
$ws.Range("R23").Value = 5484.86945658882
$ws.Range("S23").Value = 0.01203229801503209
$ws.Range("T23").Value = 0.01203229801503208
$ws.Range("G24").Value = 11.213844
$ws.Range("H24").Value = 33.641532
$ws.Range("I24").Value = 0.0601624206229761
$ws.Range("J24").Value = 0.0601624206229761
$ws.Range("M24").Value = 60.92601633333334
$ws.Range("N24").Value = 182.778049
$ws.Range("O24").Value = 0.224210932487692
$ws.Range("P24").Value = 0.224210932487692
$ws.Range("Q24").Value = 683.2148427034521
$ws.Range("R24").Value = 6148.933584331068
$ws.Range("S24").Value = 0.01348907242859422
$ws.Range("T24").Value = 0.01348907242859422
$ws.Range("G25").Value = 11.213844
$ws.Range("H25").Value = 33.641532
$ws.Range("I25").Value = 0.0601624206229761
$ws.Range("J25").Value = 0.0601624206229761
$ws.Range("M25").Value = 7.809668333333332
$ws.Range("N25").Value = 23.429005
$ws.Range("O25").Value = 0.02873998867505581
$ws.Range("P25").Value = 0.02873998867505581
$ws.Range("Q25").Value = 87.57640238173998
$ws.Range("R25").Value = 788.1876214356598
$ws.Range("S25").Value = 0.001729067287368278
$ws.Range("T25").Value = 0.001729067287368277
$ws.Range("G26").Value = 11.213844
$ws.Range("H26").Value = 33.641532
$ws.Range("I26").Value = 0.0601624206229761
$ws.Range("J26").Value = 0.0601624206229761
$ws.Range("M26").Value = 61.75795633333333
$ws.Range("N26").Value = 185.273869
$ws.Range("O26").Value = 0.2272725152794058
$ws.Range("P26").Value = 0.2272725152794058
$ws.Range("Q26").Value = 692.544088080812
$ws.Range("R26").Value = 6232.896792727307
$ws.Range("S26").Value = 0.01367326466028138
$ws.Range("T26").Value = 0.01367326466028137
